$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "42.897.45"
$ws.Range("E2").Value = "  -0.40%  "
Set-TextValue "D3" "2.300.63"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue "D5" "299.82"
$ws.Range("E5").Value = "  -0.86%  "
Set-TextValue "D6" "97.00"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -3.36%  "
Set-TextValue "D10" "35.59"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  +1.00%  "
Set-TextValue "D13" "17.85"
$ws.Range("E13").Value = "  -0.90%  "
Set-TextValue "D14" "6.76"
$ws.Range("E14").Value = "  -2.48%  "
Set-TextValue "D15" "2.656.13"
$ws.Range("E15").Value = "  -0.64%  "
Set-TextValue "D16" "2.285.08"
$ws.Range("E16").Value = "  +0.75%  "
Set-TextValue "D17" "0.777"
$ws.Range("E17").Value = "  -1.75%  "
Set-TextValue "D18" "42.824.38"
$ws.Range("E18").Value = "  -0.37%  "
Set-TextValue "D19" "12.76"
$ws.Range("E19").Value = "  -6.34%  "
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  -2.63%  "
Set-TextValue "D22" "67.80"
$ws.Range("E22").Value = "  -0.78%  "
Set-TextValue "D23" "240.03"
$ws.Range("E23").Value = "  -0.04%  "
Set-TextValue "D24" "2.13"
$ws.Range("E24").Value = "  -1.69%  "
$ws.Range("E25").Value = "  +0.07%  "
Set-TextValue "D26" "2.43"
$ws.Range("E26").Value = "  -1.28%  "
Set-TextValue "D27" "4.01"
$ws.Range("E27").Value = "  -0.48%  "
Set-TextValue "D28" "25.39"
$ws.Range("E28").Value = "  +1.75%  "
Set-TextValue "D29" "165.56"
$ws.Range("E29").Value = "  -1.85%  "
Set-TextValue "D30" "2.03"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("E31").Value = "  -1.67%  "
Set-TextValue "D32" "32.93"
$ws.Range("E32").Value = "  -1.66%  "
Set-TextValue "D33" "4.92"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  +0.08%  "
Set-TextValue "D35" "5.02"
$ws.Range("E35").Value = "  -4.08%  "
Set-TextValue "D36" "16.91"
$ws.Range("E36").Value = "  -7.99%  "
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("E42").Value = "  -1.76%  "
Set-TextValue "D43" "2.011.34"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D45" "2.16"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D46" "10.13"
$ws.Range("E46").Value = "  +0.09%  "
Set-TextValue "D47" "17.22"
$ws.Range("E47").Value = "  -2.31%  "
Set-TextValue "D48" "2.79"
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D49" "2.92"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D50" "53.43"
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D51" "2.523.69"
$ws.Range("E51").Value = "  -0.58%  "
